$wb = $excel.ActiveWorkbook

$page1 = $wb.Worksheets.Item("Page1")
$page1Data = @(
    ,@('Graduate Test Engineer / SDET', 'https://uk.indeed.com/rc/clk?jk=2316c31070f8ab52&fccid=c44218d8c2035787&vjs=3')
    ,@('SDET - Software Development Engineer in Test, Cloud (Remote)', 'https://uk.indeed.com/rc/clk?jk=24608845c54829f9&fccid=64e4cdd7435d8c42&vjs=3')
    ,@('Senior SDET', 'https://uk.indeed.com/rc/clk?jk=b258128b3092eb08&fccid=447163ca34a4b3d8&vjs=3')
    ,@('Junior SDET / Test Engineer (Python Automation)', 'https://uk.indeed.com/rc/clk?jk=f6b779cf20637696&fccid=36e692c6802aae6c&vjs=3')
    ,@('SDET', 'https://uk.indeed.com/rc/clk?jk=89a02e36c7c8699a&fccid=04f36097f3ddf22e&vjs=3')
    ,@('Software Development Engineer in Test (SDET)', 'https://uk.indeed.com/rc/clk?jk=abbf0a62614c053f&fccid=e3b893882ce68e4d&vjs=3')
    ,@('Senior Consultant - Software Development Engineer in Test (S...', 'https://uk.indeed.com/rc/clk?jk=36c06d78b4183142&fccid=9e215d88a6b33622&vjs=3')
    ,@('Senior SDET (API)', 'https://uk.indeed.com/rc/clk?jk=ec4d09e700092c5b&fccid=1d6acce86770faa4&vjs=3')
    ,@('Software Engineering in Test SDET', 'https://uk.indeed.com/rc/clk?jk=d293545f2078cefd&fccid=bfd62d072447f458&vjs=3')
    ,@('Snr SDET', 'https://uk.indeed.com/company/Qumu/jobs/Snr-Sdet-ba16e3dfbdb2d38b?fccid=77a8dc935ac1a44a&vjs=3')
    ,@('Automation Test Engineer / SDET (Java)', 'https://uk.indeed.com/rc/clk?jk=c5be219e12ad9004&fccid=c44218d8c2035787&vjs=3')
    ,@('Software Development Engineer in Test - Android (Remote, GBR...', 'https://uk.indeed.com/rc/clk?jk=bedea54b9ae2530e&fccid=64e4cdd7435d8c42&vjs=3')
    ,@('Software Development Engineer in Test', 'https://uk.indeed.com/rc/clk?jk=e65125fb4a30f046&fccid=0288b34a02030660&vjs=3')
)
for ($i = 0; $i -lt $page1Data.Count; $i++) {
    $row = $i + 2
    $page1.Cells.Item($row, 1).Value = $page1Data[$i][0]
    $page1.Cells.Item($row, 2).Value = $page1Data[$i][1]
}

$page2 = $wb.Worksheets.Item("Page2")
$page2Data = @(
    ,@('SDET - Software Development engineer in Test (Remote, GBR)', 'https://uk.indeed.com/rc/clk?jk=9db39558b28f92e1&fccid=64e4cdd7435d8c42&vjs=3')
    ,@('Software Engineer - SDET', 'https://uk.indeed.com/rc/clk?jk=82db1474bc87c15e&fccid=9526cf37acc8babe&vjs=3')
    ,@('Software Developer in Test CONTRACT', 'https://uk.indeed.com/rc/clk?jk=c04a5ccae68f8e6b&fccid=8f35589a37e69470&vjs=3')
    ,@('Senior QA Engineer C# SDET (Remote) £50k', 'https://uk.indeed.com/rc/clk?jk=2f3f5beb9b30196c&fccid=dd616958bd9ddc12&vjs=3')
    ,@('Lead Java (SDET) Software Developer In Test', 'https://uk.indeed.com/rc/clk?jk=e1b8747527ff5ea8&fccid=d6830df36521c563&vjs=3')
    ,@('Contract python test developer', 'https://uk.indeed.com/rc/clk?jk=705320533dc92838&fccid=8f35589a37e69470&vjs=3')
    ,@('C# Automation Tester', 'https://uk.indeed.com/rc/clk?jk=68ee8a1ac1c3694e&fccid=a77d4b706742e9cb&vjs=3')
    ,@('Lead Automation Tester | Remote | Up to £70,000', 'https://uk.indeed.com/rc/clk?jk=5d901ccb1ab672e1&fccid=77087bd1709a8148&vjs=3')
    ,@('Python Test Automation Engineer', 'https://uk.indeed.com/rc/clk?jk=7cf17be453a87cec&fccid=8f35589a37e69470&vjs=3')
    ,@('Senior SDET / Test Engineer (Python Automation)', 'https://uk.indeed.com/rc/clk?jk=f0eda2f03e9d44d0&fccid=36e692c6802aae6c&vjs=3')
    ,@('Full stack AWS Software Engineer', 'https://uk.indeed.com/rc/clk?jk=c4d3348ac900b88e&fccid=d50c573274bf1c1f&vjs=3')
    ,@('Lead JavaScript / React (SDET) Software Developer', 'https://uk.indeed.com/rc/clk?jk=c73c463081e061a5&fccid=d6830df36521c563&vjs=3')
    ,@('Software Development Engineer in Test - iOS (Remote, GBR)', 'https://uk.indeed.com/rc/clk?jk=4db609371d728498&fccid=64e4cdd7435d8c42&vjs=3')
)
for ($i = 0; $i -lt $page2Data.Count; $i++) {
    $row = $i + 2
    $page2.Cells.Item($row, 1).Value = $page2Data[$i][0]
    $page2.Cells.Item($row, 2).Value = $page2Data[$i][1]
}
